# Split the "	// code" run in the WhileLoops code sample into a
# "	" run (unchanged) followed by a new "// body" run, matching the
# commit that replaces the placeholder comment with real body text.

$p = $ppt.ActivePresentation

$needle = "`t// code"
$target = $null

for ($si = 1; $si -le $p.Slides.Count -and $target -eq $null; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count -and $target -eq $null; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text.Contains($needle)) {
                $target = $tr
            }
        }
    }
}

$fullText = $target.Text
$tabIdx = $fullText.IndexOf($needle)
# "// code" (without the leading tab) starts right after the tab char;
# TextRange.Characters() uses 1-based indexing, so add 2 (1 to move past
# the tab, 1 to convert from 0-based to 1-based).
$codeStart = $tabIdx + 2
$codeRange = $target.Characters($codeStart, 7)
$codeRange.Text = "// body"
